$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "modify year option to all forms": insert a new "Year of Competition"
# column right after the application-number column (new column B), shifting
# every following column one to the right.
$ws.Columns.Item(2).Insert()

# New shared-string header for the inserted column.
$ws.Range("B1").Value = "比賽年份 Year of Competition"

# The column widths for the (now widened/rearranged) leading columns were
# re-tuned by hand after the insert - set them to their final widths.
$ws.Columns.Item(1).ColumnWidth = 30.285714285714285
$ws.Columns.Item(2).ColumnWidth = 28.465401785714285
$ws.Columns.Item(3).ColumnWidth = 41.102120535714285
$ws.Columns.Item(4).ColumnWidth = 23.648995535714285
$ws.Columns.Item(5).ColumnWidth = 22.918526785714285
$ws.Columns.Item(6).ColumnWidth = 22.648995535714285
$ws.Columns.Item(7).ColumnWidth = 34.918526785714285
$ws.Columns.Item(8).ColumnWidth = 15.012276785714286
$ws.Columns.Item(9).ColumnWidth = 52.555245535714285
$ws.Columns.Item(10).ColumnWidth = 43.285714285714285

# Selection moves to A2 after the edit.
$ws.Range("A2").Select()
